$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 2 values (D2, E2) used by the new D4 formula
$ws.Range("D2").Value2 = 6000000
$ws.Range("E2").Value2 = 1

# D4 becomes a formula referencing the new cells (value stays 6000000)
$ws.Range("D4").Formula = '=D2*E2'

# Update the reference temperature K4 from 293 to 293.3 (cascades through
# L4, M4, N4 and all dependent rows 10-70)
$ws.Range("K4").Value2 = 293.3

# Update the "T=E/(CV*M)-273" label text to reflect the new constant; this
# automatically updates the shared-string table and every cell referencing it
$ws.Range("E9").Value2 = 'T=E/(CV*M)-273.3'

# Update the calc-temperature formulas from "-273" to "-273.3"
$ws.Range("E10").Formula = '=B10/$H$4/D10-273.3'
$ws.Range("E11:E70").Formula = '=B11/$H$4/D11-273.3'

# Widen the newly relevant column F and move the active selection to G6
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Range("G6").Select()
